$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need an explicit text
# number format first, otherwise Excel will coerce the assigned string
# into a numeric value (losing the literal formatting, e.g. "20.80" -> 20.8).
$textCells = @("D5","D7","D8","D9","D10","D12","D13","D14","D15","D17","D18","D19","D20","D24","D25","D26","D27","D29","D30","D31","D33","D35","D36","D37","D38","D39","D40","D41","D42","D44","D45","D46","D47","D48","D49")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values cell by cell, in sheet order.
$ws.Range('D2').Value = '28.603.73'
$ws.Range('E2').Value = '  +1.58%  '
$ws.Range('D3').Value = '1.826.51'
$ws.Range('E3').Value = '  +1.19%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '316.69'
$ws.Range('E5').Value = '  +0.03%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').Value = '0.5311'
$ws.Range('E7').Value = '  -3.05%  '
$ws.Range('D8').Value = '0.3971'
$ws.Range('E8').Value = '  +4.64%  '
$ws.Range('D9').Value = '0.07761'
$ws.Range('E9').Value = '  +3.74%  '
$ws.Range('D10').Value = '42.06'
$ws.Range('E10').Value = '  -0.26%  '
$ws.Range('E11').Value = '  +2.08%  '
$ws.Range('D12').Value = '21.13'
$ws.Range('E12').Value = '  +2.71%  '
$ws.Range('D13').Value = '6.321'
$ws.Range('E13').Value = '  +1.84%  '
$ws.Range('B14').Value = 'BinanceUSD'
$ws.Range('C14').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D14').Value = '1.002'
$ws.Range('E14').Value = '  +0.17%  '
$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '7.567'
$ws.Range('E15').Value = '  +2.78%  '
$ws.Range('D16').Value = '1.857.77'
$ws.Range('E16').Value = '  +2.99%  '
$ws.Range('D17').Value = '93.28'
$ws.Range('E17').Value = '  +3.61%  '
$ws.Range('D18').Value = '0.00001089'
$ws.Range('E18').Value = '  +2.15%  '
$ws.Range('D19').Value = '0.06623'
$ws.Range('E19').Value = '  +1.50%  '
$ws.Range('D20').Value = '17.78'
$ws.Range('E20').Value = '  +1.81%  '
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('E22').Value = '  +2.48%  '
$ws.Range('D23').Value = '28.610.81'
$ws.Range('E23').Value = '  +1.51%  '
$ws.Range('D24').Value = '11.21'
$ws.Range('E24').Value = '  -0.12%  '
$ws.Range('D25').Value = '2.237'
$ws.Range('E25').Value = '  +6.88%  '
$ws.Range('D26').Value = '20.80'
$ws.Range('E26').Value = '  +1.47%  '
$ws.Range('D27').Value = '156.93'
$ws.Range('E27').Value = '  +0.40%  '
$ws.Range('B28').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C28').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D28').Value = '2.030.67'
$ws.Range('E28').Value = '  +0.94%  '
$ws.Range('B29').Value = 'LidoDAOToken'
$ws.Range('C29').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D29').Value = '2.417'
$ws.Range('E29').Value = '  +3.07%  '
$ws.Range('D30').Value = '125.27'
$ws.Range('E30').Value = '  +2.51%  '
$ws.Range('D31').Value = '1.151'
$ws.Range('E31').Value = '  +2.42%  '
$ws.Range('E32').Value = '  +0.07%  '
$ws.Range('D33').Value = '5.746'
$ws.Range('E33').Value = '  +2.88%  '
$ws.Range('E34').Value = '  +0.17%  '
$ws.Range('D35').Value = '0.07327'
$ws.Range('E35').Value = '  +5.48%  '
$ws.Range('D36').Value = '0.2271'
$ws.Range('E36').Value = '  +1.70%  '
$ws.Range('D37').Value = '0.02350'
$ws.Range('E37').Value = '  +1.84%  '
$ws.Range('D38').Value = '8.915'
$ws.Range('D39').Value = '5.208'
$ws.Range('E39').Value = '  +2.22%  '
$ws.Range('D40').Value = '11.42'
$ws.Range('E40').Value = '  +2.02%  '
$ws.Range('D41').Value = '0.6296'
$ws.Range('E41').Value = '  +1.76%  '
$ws.Range('D42').Value = '1.195'
$ws.Range('E42').Value = '  +1.47%  '
$ws.Range('E43').Value = '  +0.03%  '
$ws.Range('D44').Value = '1.401'
$ws.Range('E44').Value = '  -1.33%  '
$ws.Range('D45').Value = '13.54'
$ws.Range('E45').Value = '  +1.29%  '
$ws.Range('D46').Value = '0.5945'
$ws.Range('E46').Value = '  +3.05%  '
$ws.Range('D47').Value = '3.721'
$ws.Range('E47').Value = '  +1.00%  '
$ws.Range('D48').Value = '125.62'
$ws.Range('E48').Value = '  +0.53%  '
$ws.Range('D49').Value = '1.999'
$ws.Range('E49').Value = '  +3.82%  '
$ws.Range('E50').Value = '  +0.29%  '
$ws.Range('E51').Value = '  +2.06%  '
